$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 1302.4375
$ws_ALC.Range("J17").Value = 1302.4375
$ws_ALC.Range("L17").Value = 3907.3125
$ws_ALC.Range("N17").Value = -4243.3125

# ALC row 64
$ws_ALC.Range("H64").Value = 58505.223
$ws_ALC.Range("J64").Value = 3174.5
$ws_ALC.Range("L64").Value = 3174.5
$ws_ALC.Range("N64").Value = -3670.5

# ALC row 67
$ws_ALC.Range("H67").Value = 58505.223
$ws_ALC.Range("J67").Value = 3174.5
$ws_ALC.Range("L67").Value = 3174.5
$ws_ALC.Range("N67").Value = -4890.5

# ALC row 74
$ws_ALC.Range("H74").Value = 3165
$ws_ALC.Range("I74").Value = 2961.111
$ws_ALC.Range("J74").Value = 5000
$ws_ALC.Range("K74").Value = 2961.111
$ws_ALC.Range("L74").Value = 5000
$ws_ALC.Range("M74").Value = -2025.111
$ws_ALC.Range("N74").Value = -6872

# ALC row 77
$ws_ALC.Range("H77").Value = 3165
$ws_ALC.Range("I77").Value = 2961.111
$ws_ALC.Range("J77").Value = 5000
$ws_ALC.Range("K77").Value = 14805.555
$ws_ALC.Range("L77").Value = 25000
$ws_ALC.Range("M77").Value = -10125.555
$ws_ALC.Range("N77").Value = -34360

# ALC row 137
$ws_ALC.Range("H137").Value = 2381.9333
$ws_ALC.Range("I137").Value = 2132.7144
$ws_ALC.Range("J137").Value = 2600
$ws_ALC.Range("K137").Value = 6398.1432
$ws_ALC.Range("L137").Value = 7800
$ws_ALC.Range("M137").Value = -3848.1432
$ws_ALC.Range("N137").Value = -12900

# ARM row 6
$ws_ARM.Range("H6").Value = 17753.637
$ws_ARM.Range("I6").Value = 37899
$ws_ARM.Range("J6").Value = 6242
$ws_ARM.Range("K6").Value = 37899
$ws_ARM.Range("L6").Value = 6242
$ws_ARM.Range("M6").Value = -37726
$ws_ARM.Range("N6").Value = -6588

# ARM row 32
$ws_ARM.Range("H32").Value = 25465.988
$ws_ARM.Range("I32").Value = 7669.9375
$ws_ARM.Range("J32").Value = 310202.8
$ws_ARM.Range("K32").Value = 7669.9375
$ws_ARM.Range("L32").Value = 310202.8
$ws_ARM.Range("M32").Value = -7382.9375
$ws_ARM.Range("N32").Value = -310776.8

# ARM row 74
$ws_ARM.Range("H74").Value = 1076.8182
$ws_ARM.Range("I74").Value = 926.4
$ws_ARM.Range("J74").Value = 1202.1666
$ws_ARM.Range("K74").Value = 926.4
$ws_ARM.Range("L74").Value = 1202.1666
$ws_ARM.Range("M74").Value = -52.39999999999998
$ws_ARM.Range("N74").Value = -2950.1666

# ARM row 77
$ws_ARM.Range("H77").Value = 1076.8182
$ws_ARM.Range("I77").Value = 926.4
$ws_ARM.Range("J77").Value = 1202.1666
$ws_ARM.Range("K77").Value = 4632
$ws_ARM.Range("L77").Value = 6010.833000000001
$ws_ARM.Range("M77").Value = -264
$ws_ARM.Range("N77").Value = -14746.833

# ARM row 80
$ws_ARM.Range("H80").Value = 24456.363
$ws_ARM.Range("J80").Value = 24902
$ws_ARM.Range("L80").Value = 24902
$ws_ARM.Range("N80").Value = -26898

# ARM row 83
$ws_ARM.Range("H83").Value = 24456.363
$ws_ARM.Range("J83").Value = 24902
$ws_ARM.Range("L83").Value = 74706
$ws_ARM.Range("N83").Value = -84690

# ARM row 102
$ws_ARM.Range("H102").Value = 65796.31
$ws_ARM.Range("I102").Value = 126984.875
$ws_ARM.Range("J102").Value = 4607.75
$ws_ARM.Range("K102").Value = 126984.875
$ws_ARM.Range("L102").Value = 4607.75
$ws_ARM.Range("M102").Value = -125362.875
$ws_ARM.Range("N102").Value = -7851.75

# ARM row 122
$ws_ARM.Range("H122").Value = 1996.0588
$ws_ARM.Range("I122").Value = 1687.6086
$ws_ARM.Range("K122").Value = 5062.825800000001
$ws_ARM.Range("M122").Value = -2612.825800000001

# ARM row 132
$ws_ARM.Range("H132").Value = 12542.444
$ws_ARM.Range("I132").Value = 15030.581
$ws_ARM.Range("J132").Value = 2816.0908
$ws_ARM.Range("K132").Value = 45091.743
$ws_ARM.Range("L132").Value = 8448.2724
$ws_ARM.Range("M132").Value = -42561.743
$ws_ARM.Range("N132").Value = -13508.2724

# BSM row 92
$ws_BSM.Range("H92").Value = 15000
$ws_BSM.Range("J92").Value = 15000
$ws_BSM.Range("L92").Value = 15000
$ws_BSM.Range("N92").Value = -19992

# BSM row 99
$ws_BSM.Range("H99").Value = 2114.5757
$ws_BSM.Range("I99").Value = 1232.2222
$ws_BSM.Range("J99").Value = 2445.4583
$ws_BSM.Range("K99").Value = 1232.2222
$ws_BSM.Range("L99").Value = 2445.4583
$ws_BSM.Range("M99").Value = 265.7778000000001
$ws_BSM.Range("N99").Value = -5441.4583

# BSM row 105
$ws_BSM.Range("H105").Value = 85060.75
$ws_BSM.Range("I105").Value = 64236.688
$ws_BSM.Range("J105").Value = 126708.875
$ws_BSM.Range("K105").Value = 64236.688
$ws_BSM.Range("L105").Value = 126708.875
$ws_BSM.Range("M105").Value = -62489.688
$ws_BSM.Range("N105").Value = -130202.875

# BSM row 134
$ws_BSM.Range("H134").Value = 2788.9148
$ws_BSM.Range("I134").Value = 2738.2104
$ws_BSM.Range("J134").Value = 3003
$ws_BSM.Range("K134").Value = 8214.6312
$ws_BSM.Range("L134").Value = 9009
$ws_BSM.Range("M134").Value = -5679.6312
$ws_BSM.Range("N134").Value = -14079

# CRP row 12
$ws_CRP.Range("H12").Value = 5168.6665
$ws_CRP.Range("I12").Value = 506
$ws_CRP.Range("J12").Value = 7500
$ws_CRP.Range("K12").Value = 506
$ws_CRP.Range("L12").Value = 7500
$ws_CRP.Range("M12").Value = -336
$ws_CRP.Range("N12").Value = -7840

# CRP row 31
$ws_CRP.Range("H31").Value = 33220.586
$ws_CRP.Range("I31").Value = 966.9048
$ws_CRP.Range("J31").Value = 60313.68
$ws_CRP.Range("K31").Value = 966.9048
$ws_CRP.Range("L31").Value = 60313.68
$ws_CRP.Range("M31").Value = -671.9048
$ws_CRP.Range("N31").Value = -60903.68

# CRP row 34
$ws_CRP.Range("H34").Value = 33220.586
$ws_CRP.Range("I34").Value = 966.9048
$ws_CRP.Range("J34").Value = 60313.68
$ws_CRP.Range("K34").Value = 966.9048
$ws_CRP.Range("L34").Value = 60313.68
$ws_CRP.Range("M34").Value = -764.9048
$ws_CRP.Range("N34").Value = -60717.68

# CRP row 105
$ws_CRP.Range("H105").Value = 1077.7059
$ws_CRP.Range("I105").Value = 1094.3
$ws_CRP.Range("J105").Value = 1054
$ws_CRP.Range("K105").Value = 1094.3
$ws_CRP.Range("L105").Value = 1054
$ws_CRP.Range("M105").Value = 652.7
$ws_CRP.Range("N105").Value = -4548

# CUL row 37
$ws_CUL.Range("H37").Value = 610829.25
$ws_CUL.Range("J37").Value = 610829.25
$ws_CUL.Range("L37").Value = 1832487.75
$ws_CUL.Range("N37").Value = -1832711.75

# CUL row 92
$ws_CUL.Range("H92").Value = 495.75
$ws_CUL.Range("I92").Value = 490
$ws_CUL.Range("J92").Value = 501.5
$ws_CUL.Range("K92").Value = 1470
$ws_CUL.Range("L92").Value = 1504.5
$ws_CUL.Range("M92").Value = -222
$ws_CUL.Range("N92").Value = -4000.5

# CUL row 113
$ws_CUL.Range("H113").Value = 976.2692
$ws_CUL.Range("I113").Value = 1416.1818
$ws_CUL.Range("J113").Value = 653.6667
$ws_CUL.Range("K113").Value = 4248.5454
$ws_CUL.Range("L113").Value = 1961.0001
$ws_CUL.Range("M113").Value = -2078.5454
$ws_CUL.Range("N113").Value = -6301.0001

# CUL row 131
$ws_CUL.Range("H131").Value = 618149.5
$ws_CUL.Range("I131").Value = 670
$ws_CUL.Range("J131").Value = 694277.1
$ws_CUL.Range("K131").Value = 2010
$ws_CUL.Range("L131").Value = 2082831.3
$ws_CUL.Range("M131").Value = 3030
$ws_CUL.Range("N131").Value = -2092911.3

# CUL row 133
$ws_CUL.Range("H133").Value = 3759.5833
$ws_CUL.Range("I133").Value = 1280
$ws_CUL.Range("J133").Value = 4999.375
$ws_CUL.Range("K133").Value = 3840
$ws_CUL.Range("L133").Value = 14998.125
$ws_CUL.Range("M133").Value = 1220
$ws_CUL.Range("N133").Value = -25118.125

# CUL row 134
$ws_CUL.Range("H134").Value = 3889.318
$ws_CUL.Range("I134").Value = 2168.5334
$ws_CUL.Range("J134").Value = 7576.7144
$ws_CUL.Range("K134").Value = 6505.600199999999
$ws_CUL.Range("L134").Value = 22730.1432
$ws_CUL.Range("M134").Value = -1435.600199999999
$ws_CUL.Range("N134").Value = -32870.1432

# CUL row 136
$ws_CUL.Range("H136").Value = 1370.3334
$ws_CUL.Range("I136").Value = 1641
$ws_CUL.Range("J136").Value = 829
$ws_CUL.Range("K136").Value = 4923
$ws_CUL.Range("L136").Value = 2487
$ws_CUL.Range("M136").Value = 177
$ws_CUL.Range("N136").Value = -12687

# CUL row 137
$ws_CUL.Range("H137").Value = 3746980.8
$ws_CUL.Range("I137").Value = 93137.27
$ws_CUL.Range("J137").Value = 6258998
$ws_CUL.Range("K137").Value = 279411.81
$ws_CUL.Range("L137").Value = 18776994
$ws_CUL.Range("M137").Value = -274311.81
$ws_CUL.Range("N137").Value = -18787194

# CUL row 138
$ws_CUL.Range("H138").Value = 8465.666999999999
$ws_CUL.Range("I138").Value = 14123.75
$ws_CUL.Range("J138").Value = 1999.2858
$ws_CUL.Range("K138").Value = 42371.25
$ws_CUL.Range("L138").Value = 5997.857400000001
$ws_CUL.Range("M138").Value = -37231.25
$ws_CUL.Range("N138").Value = -16277.8574

# CUL row 140
$ws_CUL.Range("H140").Value = 5190.5356
$ws_CUL.Range("I140").Value = 6367
$ws_CUL.Range("J140").Value = 2249.375
$ws_CUL.Range("K140").Value = 19101
$ws_CUL.Range("L140").Value = 6748.125
$ws_CUL.Range("M140").Value = -13921
$ws_CUL.Range("N140").Value = -17108.125

# GSM row 132
$ws_GSM.Range("H132").Value = 2620.9512
$ws_GSM.Range("I132").Value = 1896.7097
$ws_GSM.Range("J132").Value = 4866.1
$ws_GSM.Range("K132").Value = 5690.1291
$ws_GSM.Range("L132").Value = 14598.3
$ws_GSM.Range("M132").Value = -3160.1291
$ws_GSM.Range("N132").Value = -19658.3

# LTW row 7
$ws_LTW.Range("H7").Value = 2795.111
$ws_LTW.Range("I7").Value = 1791.2
$ws_LTW.Range("J7").Value = 4050
$ws_LTW.Range("K7").Value = 1791.2
$ws_LTW.Range("L7").Value = 4050
$ws_LTW.Range("M7").Value = -1679.2
$ws_LTW.Range("N7").Value = -4274

# LTW row 82
$ws_LTW.Range("H82").Value = 2154.7144
$ws_LTW.Range("I82").Value = 1288.3334
$ws_LTW.Range("J82").Value = 2804.5
$ws_LTW.Range("K82").Value = 1288.3334
$ws_LTW.Range("L82").Value = 2804.5
$ws_LTW.Range("M82").Value = -927.3334
$ws_LTW.Range("N82").Value = -3526.5

# LTW row 85
$ws_LTW.Range("H85").Value = 2154.7144
$ws_LTW.Range("I85").Value = 1288.3334
$ws_LTW.Range("J85").Value = 2804.5
$ws_LTW.Range("K85").Value = 1288.3334
$ws_LTW.Range("L85").Value = 2804.5
$ws_LTW.Range("M85").Value = -40.33339999999998
$ws_LTW.Range("N85").Value = -5300.5

# LTW row 126
$ws_LTW.Range("H126").Value = 2795.111
$ws_LTW.Range("I126").Value = 1791.2
$ws_LTW.Range("J126").Value = 4050
$ws_LTW.Range("K126").Value = 5373.6
$ws_LTW.Range("L126").Value = 12150
$ws_LTW.Range("M126").Value = -2903.6
$ws_LTW.Range("N126").Value = -17090

# WVR row 123
$ws_WVR.Range("H123").Value = 19571.428
$ws_WVR.Range("J123").Value = 19571.428
$ws_WVR.Range("L123").Value = 19571.428
$ws_WVR.Range("N123").Value = -29371.428

# WVR row 126
$ws_WVR.Range("H126").Value = 1551.76
$ws_WVR.Range("I126").Value = 1364.95
$ws_WVR.Range("J126").Value = 2299
$ws_WVR.Range("K126").Value = 4094.85
$ws_WVR.Range("L126").Value = 6897
$ws_WVR.Range("M126").Value = -1624.85
$ws_WVR.Range("N126").Value = -11837

# WVR row 132
$ws_WVR.Range("H132").Value = 4992.5366
$ws_WVR.Range("I132").Value = 2800.077
$ws_WVR.Range("J132").Value = 8792.799999999999
$ws_WVR.Range("K132").Value = 8400.231
$ws_WVR.Range("L132").Value = 26378.4
$ws_WVR.Range("M132").Value = -5870.231
$ws_WVR.Range("N132").Value = -31438.4
